$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-based values for rows 2-7 (columns A:T)
# Column order: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T numeric edge-weight columns

$data = @(
    @("ECs","Col9a2","Mag","MuSCs",2,0.6666666666666666,0.1981143333333334,0.5943430000000001,0.1714078479099734,0.1714078479099734,3,1,0.8180823333333334,2.454247,0.5115352725808422,0.5115352725808422,0.1620738360801111,1.458664524721,0.0876811602031238,0.0876811602031238),
    @("ECs","Col9a2","Mag","Resolving-Mac",2,0.6666666666666666,0.1981143333333334,0.5943430000000001,0.1714078479099734,0.1714078479099734,3,1,0.7811863333333333,2.343559,0.4884647274191578,0.4884647274191579,0.1547642096374445,1.392877886737,0.08372668770684963,0.08372668770684964),
    @("FAPs","Col9a2","Mag","MuSCs",3,1,0.8244333333333334,2.4733,0.7132969181697054,0.7132969181697054,3,1,0.8180823333333334,2.454247,0.5115352725808422,0.5115352725808422,0.6744543450111111,6.0700891051,0.3648765334670149,0.3648765334670149),
    @("FAPs","Col9a2","Mag","Resolving-Mac",3,1,0.8244333333333334,2.4733,0.7132969181697054,0.7132969181697054,3,1,0.7811863333333333,2.343559,0.4884647274191578,0.4884647274191579,0.6440360527444444,5.7963244747,0.3484203847026905,0.3484203847026905),
    @("MuSCs","Col9a2","Mag","MuSCs",2,0.6666666666666666,0.133259,0.399777,0.1152952339203211,0.1152952339203211,3,1,0.8180823333333334,2.454247,0.5115352725808422,0.5115352725808422,0.1090168336576667,0.981151502919,0.05897757891070344,0.05897757891070344),
    @("MuSCs","Col9a2","Mag","Resolving-Mac",2,0.6666666666666666,0.133259,0.399777,0.1152952339203211,0.1152952339203211,3,1,0.7811863333333333,2.343559,0.4884647274191578,0.4884647274191579,0.1041001095936667,0.936900986343,0.05631765500961771,0.05631765500961771)
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $addr = "$($cols[$j])$rowNum"
        $ws.Range($addr).Value = $rowVals[$j]
    }
}

$wb.Save()
